# Auto-generated Excel COM-interop edit script
# Applies numeric cell-value corrections (Sophia_Profits workbook refresh)
$wb = $excel.ActiveWorkbook

# --- ALC!row 17 (hunk 0) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4200
$ws.Range("J17").Value = 4200
$ws.Range("L17").Value = 12600
$ws.Range("N17").Value = -12936

# --- ALC!row 28 (hunk 1) ---
$ws.Range("H28").Value = 412.75
$ws.Range("I28").Value = 215
$ws.Range("K28").Value = 215
$ws.Range("M28").Value = 270

# --- ALC!row 88 (hunk 2) ---
$ws.Range("H88").Value = 1699
$ws.Range("J88").Value = 1833.3334
$ws.Range("L88").Value = 1833.3334
$ws.Range("N88").Value = -2645.3334

# --- ALC!row 91 (hunk 3) ---
$ws.Range("H91").Value = 1699
$ws.Range("J91").Value = 1833.3334
$ws.Range("L91").Value = 1833.3334
$ws.Range("N91").Value = -4641.3334

# --- ALC!row 98 (hunk 4) ---
$ws.Range("H98").Value = 6224.25
$ws.Range("I98").Value = 6224.25
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 6224.25
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -4726.25
$ws.Range("N98").ClearContents()

# --- ALC!row 107 (hunk 5) ---
$ws.Range("H107").Value = 7109
$ws.Range("I107").Value = 6386.25
$ws.Range("K107").Value = 6386.25
$ws.Range("M107").Value = -4466.25

# --- ALC!row 111 (hunk 6) ---
$ws.Range("H111").Value = 411.33334
$ws.Range("I111").Value = 411.33334
$ws.Range("K111").Value = 1234.00002
$ws.Range("M111").Value = 1832.99998

# --- ALC!row 122 (hunk 7) ---
$ws.Range("H122").Value = 6224.25
$ws.Range("I122").Value = 6224.25
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 18672.75
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -16222.75
$ws.Range("N122").ClearContents()

# --- ARM!row 2 (hunk 8) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4053.25
$ws.Range("I2").Value = 6250
$ws.Range("J2").Value = 1856.5
$ws.Range("K2").Value = 6250
$ws.Range("L2").Value = 1856.5
$ws.Range("M2").Value = -6137
$ws.Range("N2").Value = -2082.5

# --- ARM!row 32 (hunk 9) ---
$ws.Range("H32").Value = 6906.364
$ws.Range("I32").Value = 6906.364
$ws.Range("K32").Value = 6906.364
$ws.Range("M32").Value = -6619.364

# --- ARM!row 116 (hunk 10) ---
$ws.Range("H116").Value = 4053.25
$ws.Range("I116").Value = 6250
$ws.Range("J116").Value = 1856.5
$ws.Range("K116").Value = 6250
$ws.Range("L116").Value = 1856.5
$ws.Range("M116").Value = -3956
$ws.Range("N116").Value = -6444.5

# --- BSM!row 3 (hunk 11) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4053.25
$ws.Range("I3").Value = 6250
$ws.Range("J3").Value = 1856.5
$ws.Range("K3").Value = 6250
$ws.Range("L3").Value = 1856.5
$ws.Range("M3").Value = -6136
$ws.Range("N3").Value = -2084.5

# --- BSM!row 64 (hunk 12) ---
$ws.Range("H64").Value = 3469.6
$ws.Range("I64").Value = 574
$ws.Range("J64").Value = 5400
$ws.Range("K64").Value = 574
$ws.Range("L64").Value = 5400
$ws.Range("M64").Value = -349
$ws.Range("N64").Value = -5850

# --- BSM!row 67 (hunk 13) ---
$ws.Range("H67").Value = 3469.6
$ws.Range("I67").Value = 574
$ws.Range("J67").Value = 5400
$ws.Range("K67").Value = 574
$ws.Range("L67").Value = 5400
$ws.Range("M67").Value = 206
$ws.Range("N67").Value = -6960

# --- BSM!row 86 (hunk 14) ---
$ws.Range("H86").Value = 4814
$ws.Range("I86").Value = 2891.5
$ws.Range("J86").Value = 6095.6665
$ws.Range("K86").Value = 2891.5
$ws.Range("L86").Value = 6095.6665
$ws.Range("M86").Value = -1768.5
$ws.Range("N86").Value = -8341.666499999999

# --- BSM!row 89 (hunk 15) ---
$ws.Range("H89").Value = 4814
$ws.Range("I89").Value = 2891.5
$ws.Range("J89").Value = 6095.6665
$ws.Range("K89").Value = 14457.5
$ws.Range("L89").Value = 30478.3325
$ws.Range("M89").Value = -8841.5
$ws.Range("N89").Value = -41710.3325

# --- BSM!row 94 (hunk 16) ---
$ws.Range("H94").Value = 2478
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()

# --- CRP!row 16 (hunk 17) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()

# --- CRP!row 31 (hunk 18) ---
$ws.Range("H31").Value = 2990.75
$ws.Range("I31").Value = 2787.8333
$ws.Range("K31").Value = 2787.8333
$ws.Range("M31").Value = -2492.8333

# --- CRP!row 34 (hunk 19) ---
$ws.Range("H34").Value = 2990.75
$ws.Range("I34").Value = 2787.8333
$ws.Range("K34").Value = 2787.8333
$ws.Range("M34").Value = -2585.8333

# --- CRP!row 113 (hunk 20) ---
$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()

# --- CRP!row 134 (hunk 21) ---
$ws.Range("H134").Value = 6796.6665
$ws.Range("J134").Value = 2932.6667
$ws.Range("L134").Value = 8798.000100000001
$ws.Range("N134").Value = -13868.0001

# --- CUL!row 13 (hunk 22) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 1040
$ws.Range("I13").Value = 60
$ws.Range("J13").Value = 3000
$ws.Range("K13").Value = 180
$ws.Range("L13").Value = 9000
$ws.Range("M13").Value = -12
$ws.Range("N13").Value = -9336

# --- LTW!row 22 (hunk 23) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 6033.3687
$ws.Range("I22").Value = 3989.4546
$ws.Range("J22").Value = 8843.75
$ws.Range("K22").Value = 3989.4546
$ws.Range("L22").Value = 8843.75
$ws.Range("M22").Value = -3694.4546
$ws.Range("N22").Value = -9433.75

# --- LTW!row 27 (hunk 24) ---
$ws.Range("H27").Value = 6033.3687
$ws.Range("I27").Value = 3989.4546
$ws.Range("J27").Value = 8843.75
$ws.Range("K27").Value = 3989.4546
$ws.Range("L27").Value = 8843.75
$ws.Range("M27").Value = -3882.4546
$ws.Range("N27").Value = -9057.75

# --- LTW!row 46 (hunk 25) ---
$ws.Range("H46").Value = 8539.923000000001
$ws.Range("J46").Value = 31153
$ws.Range("L46").Value = 31153
$ws.Range("N46").Value = -31529

# --- LTW!row 55 (hunk 26) ---
$ws.Range("H55").Value = 5200
$ws.Range("I55").Value = 5625
$ws.Range("J55").Value = 3500
$ws.Range("K55").Value = 5625
$ws.Range("L55").Value = 3500
$ws.Range("M55").Value = -5452
$ws.Range("N55").Value = -3846

# --- WVR!row 69 (hunk 27) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 18722
$ws.Range("J69").Value = 18722
$ws.Range("L69").Value = 18722
$ws.Range("N69").Value = -20220

# --- WVR!row 72 (hunk 28) ---
$ws.Range("H72").Value = 18722
$ws.Range("J72").Value = 18722
$ws.Range("L72").Value = 56166
$ws.Range("N72").Value = -63654

# --- WVR!row 107 (hunk 29) ---
$ws.Range("H107").Value = 996.7778
$ws.Range("I107").Value = 828.6667
$ws.Range("J107").Value = 1333
$ws.Range("K107").Value = 2486.0001
$ws.Range("L107").Value = 3999
$ws.Range("M107").Value = -566.0001000000002
$ws.Range("N107").Value = -7839

# --- WVR!row 126 (hunk 30) ---
$ws.Range("H126").Value = 1004.5
$ws.Range("J126").Value = 1004.5
$ws.Range("L126").Value = 3013.5
$ws.Range("N126").Value = -7953.5

Write-Host "Applied 31 row updates across ALC/ARM/BSM/CRP/CUL/LTW/WVR"
